$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells we are about to touch to stay as plain
# text, since several values (e.g. "6.41", "0.999", "0.0000228") would
# otherwise be auto-converted into numbers by Excel, losing the exact
# textual representation used in the source data (thousand-dot groups,
# fixed decimal places, etc).
$dCells = @("D2", "D3", "D5", "D6", "D7", "D8", "D9", "D10", "D12", "D13", "D14", "D15", "D16", "D17", "D19", "D20", "D22", "D23", "D24", "D25", "D27", "D28", "D31", "D32", "D33", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D50", "D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '64.676.32'
$ws.Range("E2").Value = '  +3.83%  '
$ws.Range("D3").Value = '3.086.92'
$ws.Range("E3").Value = '  +1.91%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '559.16'
$ws.Range("E5").Value = '  +3.25%  '
$ws.Range("D6").Value = '143.68'
$ws.Range("E6").Value = '  +7.45%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '3.084.94'
$ws.Range("E8").Value = '  +2.10%  '
$ws.Range("D9").Value = '0.499'
$ws.Range("E9").Value = '  +1.05%  '
$ws.Range("D10").Value = '6.41'
$ws.Range("E10").Value = '  +4.45%  '
$ws.Range("E11").Value = '  +2.90%  '
$ws.Range("D12").Value = '0.468'
$ws.Range("E12").Value = '  +5.00%  '
$ws.Range("D13").Value = '0.0000228'
$ws.Range("E13").Value = '  +2.63%  '
$ws.Range("D14").Value = '35.18'
$ws.Range("D15").Value = '3.602.05'
$ws.Range("E15").Value = '  +2.47%  '
$ws.Range("D16").Value = '64.647.03'
$ws.Range("E16").Value = '  +3.82%  '
$ws.Range("D17").Value = '3.089.75'
$ws.Range("E17").Value = '  +2.20%  '
$ws.Range("E18").Value = '  +0.91%  '
$ws.Range("D19").Value = '6.75'
$ws.Range("E19").Value = '  +1.90%  '
$ws.Range("D20").Value = '477.83'
$ws.Range("E21").Value = '  +3.40%  '
$ws.Range("D22").Value = '0.683'
$ws.Range("E22").Value = '  +1.47%  '
$ws.Range("D23").Value = '7.57'
$ws.Range("E23").Value = '  +7.26%  '
$ws.Range("D24").Value = '13.40'
$ws.Range("E24").Value = '  +10.68%  '
$ws.Range("D25").Value = '81.08'
$ws.Range("E25").Value = '  +0.30%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").Value = '2.77'
$ws.Range("E27").Value = '  +2.47%  '
$ws.Range("D28").Value = '8.17'
$ws.Range("E28").Value = '  +5.38%  '
$ws.Range("E29").Value = '  +5.98%  '
$ws.Range("E30").Value = '  +0.10%  '
$ws.Range("D31").Value = '26.03'
$ws.Range("E31").Value = '  +1.32%  '
$ws.Range("D32").Value = '1.15'
$ws.Range("E32").Value = '  +2.45%  '
$ws.Range("D33").Value = '2.46'
$ws.Range("E33").Value = '  +4.16%  '
$ws.Range("E34").Value = '  -1.03%  '
$ws.Range("E35").Value = '  +4.50%  '
$ws.Range("D36").Value = '54.94'
$ws.Range("E36").Value = '  -0.17%  '
$ws.Range("D37").Value = '460.66'
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '0.0830'
$ws.Range("E38").Value = '  +3.58%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.0406'
$ws.Range("E39").Value = '  +4.77%  '
$ws.Range("D40").Value = '2.93'
$ws.Range("E40").Value = '  +18.44%  '
$ws.Range("D41").Value = '2.972.03'
$ws.Range("E41").Value = '  -6.14%  '
$ws.Range("D42").Value = '8.22'
$ws.Range("E42").Value = '  +1.60%  '
$ws.Range("E43").Value = '  -2.43%  '
$ws.Range("D44").Value = '27.95'
$ws.Range("E44").Value = '  +5.72%  '
$ws.Range("D45").Value = '0.259'
$ws.Range("E45").Value = '  +5.63%  '
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("E47").Value = '  +7.97%  '
$ws.Range("E48").Value = '  +2.83%  '
$ws.Range("E49").Value = '  +4.42%  '
$ws.Range("D50").Value = '116.85'
$ws.Range("E50").Value = '  +2.48%  '
$ws.Range("D51").Value = '2.05'
$ws.Range("E51").Value = '  +1.65%  '
